# Fix Master Data Validation
# Reworks the Supplier Import Template header row:
#   - existing required columns get a "*" suffix (SupplierCode*, SupplierName*, Address*)
#   - ContactPhone / ContactEmail stay put
#   - five new columns are appended: DeliveryTerm, PaymentTerm, CurrencyCode*, TaxRate*, CountryCode*
#   - TaxRate* gets a numeric "0.00" format
#   - sheet view zoom + selection + column widths are refreshed to match the wider sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row text -------------------------------------------------
# Write the brand-new columns (F, G) before touching A:C so the shared
# string table ends up ordered the same way the source workbook has it
# (surviving strings first, then new ones in first-write order).
$ws.Range("F1").Value = "DeliveryTerm"
$ws.Range("G1").Value = "PaymentTerm"

$ws.Range("A1").Value = "SupplierCode*"
$ws.Range("B1").Value = "SupplierName*"
$ws.Range("C1").Value = "Address*"
# D1 (ContactPhone) / E1 (ContactEmail) are already correct - leave as-is.

$ws.Range("H1").Value = "CurrencyCode*"
$ws.Range("I1").Value = "TaxRate*"
$ws.Range("I1").NumberFormat = "0.00"
$ws.Range("J1").Value = "CountryCode*"

# --- Column widths (OOXML "width" = COM ColumnWidth + 5/6) -----------
$ws.Columns.Item(1).ColumnWidth  = (12.6640625  - 5/6)
$ws.Columns.Item(4).ColumnWidth  = (12.77734375 - 5/6)
$ws.Columns.Item(5).ColumnWidth  = (12          - 5/6)
$ws.Columns.Item(6).ColumnWidth  = (11.5546875  - 5/6)
$ws.Columns.Item(7).ColumnWidth  = (12.109375   - 5/6)
$ws.Columns.Item(8).ColumnWidth  = (13.88671875 - 5/6)
$ws.Columns.Item(9).ColumnWidth  = (8.33203125  - 5/6)
$ws.Columns.Item(10).ColumnWidth = (12.88671875 - 5/6)

# --- View: zoom + active selection ------------------------------------
$excel.ActiveWindow.Zoom = 145
$ws.Range("F4").Select()
